$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.011.98"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "2.635.71"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.80"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").Value = "'156.18"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.651"
$ws.Range("E7").Value = "  +3.94%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.122"
$ws.Range("E9").Value = "  -3.79%  "
$ws.Range("D10").Value = "'5.82"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'0.388"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("D12").Value = "'0.156"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'28.70"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "3.115.50"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'0.0000186"
$ws.Range("E15").Value = "  -5.61%  "
$ws.Range("D16").Value = "63.899.10"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").Value = "2.629.59"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "'12.22"
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("D19").Value = "'4.68"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "'7.61"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "'348.10"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "'67.80"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("E24").Value = "  +6.13%  "
$ws.Range("D25").Value = "'0.0000110"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").Value = "'9.38"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("D27").Value = "'578.73"
$ws.Range("E27").Value = "  +9.27%  "
$ws.Range("D28").Value = "'1.59"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'2.09"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'1.72"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "'6.51"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'5.32"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "'0.412"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "'20.04"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "'152.31"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'41.95"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'159.55"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'2.39"
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("D45").Value = "'4.00"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").Value = "'23.25"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").Value = "'0.0599"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("D49").Value = "'0.634"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'0.0254"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "'19.20"
$ws.Range("E51").Value = "  -2.66%  "
